# Update the title-slide author byline:
#   "PRESENTED BY : JENIFER.C"     -> "PRESENTED BY : Alan .G"
#   "REGISTER NO: 312207857"       -> "REGISTER NO: 312207882"
#
# The replacement is done by re-setting only the trailing characters of
# each paragraph (via TextRange.Characters), which is what PowerPoint
# itself does when a user selects a sub-string and retypes it - it keeps
# the untouched leading run intact and creates a new run for the
# replaced characters (inheriting the original run's formatting).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the subtitle placeholder that holds the byline text (normally
# Shape 2 / "Subtitle 2"), falling back to a text search if the deck
# layout ever changes.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.HasTextFrame) {
        $t = $cand.TextFrame.TextRange.Text
        if ($t -like "*PRESENTED BY*") {
            $target = $cand
            break
        }
    }
}
if ($target -eq $null) {
    $target = $s.Shapes.Item(2)
}

$tr = $target.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count

for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $text = $para.Text

    if ($text -like "*PRESENTED BY*JENIFER.C*") {
        $oldSuffix = " JENIFER.C"
        $newSuffix = " Alan .G"
        $idx = $text.IndexOf($oldSuffix)
        if ($idx -ge 0) {
            $sub = $para.Characters($idx + 1, $oldSuffix.Length)
            $sub.Text = $newSuffix
        }
    }
    elseif ($text -like "*REGISTER NO*312207857*") {
        $oldSuffix = "857"
        $newSuffix = "882"
        $idx = $text.LastIndexOf($oldSuffix)
        if ($idx -ge 0) {
            $sub = $para.Characters($idx + 1, $oldSuffix.Length)
            $sub.Text = $newSuffix
        }
    }
}
